$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''67.269.09'
$ws.Range('E2').Value = '''  -3.32%  '

$ws.Range('D3').Value = '''3.491.21'
$ws.Range('E3').Value = '''  -4.90%  '

$ws.Range('E4').Value = '''  +0.08%  '

$ws.Range('D5').Value = '''605.05'
$ws.Range('E5').Value = '''  -2.30%  '

$ws.Range('D6').Value = '''148.58'
$ws.Range('E6').Value = '''  -6.87%  '

$ws.Range('D7').Value = '''3.491.45'
$ws.Range('E7').Value = '''  -4.79%  '

$ws.Range('E8').Value = '''  -0.08%  '

$ws.Range('E9').Value = '''  -3.06%  '

$ws.Range('E10').Value = '''  -3.89%  '

$ws.Range('E11').Value = '''  -2.89%  '

$ws.Range('E12').Value = '''  -4.02%  '

$ws.Range('E13').Value = '''  -4.79%  '

$ws.Range('D14').Value = '''4.080.97'
$ws.Range('E14').Value = '''  -4.87%  '

$ws.Range('D15').Value = '''31.40'
$ws.Range('E15').Value = '''  -3.01%  '

$ws.Range('D16').Value = '''3.491.51'
$ws.Range('E16').Value = '''  -4.04%  '

$ws.Range('D17').Value = '''67.181.57'
$ws.Range('E17').Value = '''  -3.53%  '

$ws.Range('E19').Value = '''  -1.89%  '

$ws.Range('D20').Value = '''15.02'
$ws.Range('E20').Value = '''  -5.21%  '

$ws.Range('D21').Value = '''445.99'
$ws.Range('E21').Value = '''  -5.18%  '

$ws.Range('E22').Value = '''  -12.57%  '

$ws.Range('D23').Value = '''0.620'
$ws.Range('E23').Value = '''  -4.26%  '

$ws.Range('D24').Value = '''77.08'
$ws.Range('E24').Value = '''  -3.43%  '

$ws.Range('E25').Value = '''  +3.68%  '

$ws.Range('E26').Value = '''  +0.10%  '

$ws.Range('D27').Value = '''3.629.77'
$ws.Range('E27').Value = '''  -4.92%  '

$ws.Range('D28').Value = '''10.12'
$ws.Range('E28').Value = '''  -8.35%  '

$ws.Range('D29').Value = '''8.27'
$ws.Range('E29').Value = '''  -5.06%  '

$ws.Range('D30').Value = '''2.47'
$ws.Range('E30').Value = '''  -4.41%  '

$ws.Range('E31').Value = '''  -6.36%  '

$ws.Range('D32').Value = '''1.00'
$ws.Range('E32').Value = '''  +0.04%  '

$ws.Range('E33').Value = '''  -0.23%  '

$ws.Range('D34').Value = '''25.63'
$ws.Range('E34').Value = '''  -3.55%  '

$ws.Range('D35').Value = '''6.13'
$ws.Range('E35').Value = '''  -4.09%  '

$ws.Range('E36').Value = '''  -6.55%  '

$ws.Range('D37').Value = '''3.480.13'
$ws.Range('E37').Value = '''  -5.22%  '

$ws.Range('D38').Value = '''7.99'
$ws.Range('E38').Value = '''  -3.35%  '

$ws.Range('E39').Value = '''  +0.06%  '

$ws.Range('D40').Value = '''1.00'
$ws.Range('E40').Value = '''  +0.13%  '

$ws.Range('B41').Value = '''Monero'
$ws.Range('C41').Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').Value = '''173.83'
$ws.Range('E41').Value = '''  -2.35%  '

$ws.Range('B42').Value = '''Stacks'
$ws.Range('C42').Value = '''https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '''2.18'
$ws.Range('E42').Value = '''  -0.78%  '

$ws.Range('D43').Value = '''0.0877'
$ws.Range('E43').Value = '''  -1.39%  '

$ws.Range('E44').Value = '''  -6.72%  '

$ws.Range('E45').Value = '''  -4.64%  '

$ws.Range('D46').Value = '''45.42'
$ws.Range('E46').Value = '''  -2.66%  '

$ws.Range('D47').Value = '''27.02'
$ws.Range('E47').Value = '''  -6.03%  '

$ws.Range('D48').Value = '''1.25'
$ws.Range('E48').Value = '''  +4.13%  '

$ws.Range('E49').Value = '''  -5.67%  '

$ws.Range('D50').Value = '''7.53'
$ws.Range('E50').Value = '''  -4.03%  '

$ws.Range('E51').Value = '''  -3.20%  '
